$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell K1, matching the bold/centered style used by the
# other header cells in row 1 (e.g. J1).
$ws.Range("K1").Value = "PhylogenySorting"
$ws.Range("K1").Font.Bold = $ws.Range("J1").Font.Bold
$ws.Range("K1").HorizontalAlignment = $ws.Range("J1").HorizontalAlignment

# Add new cell K3 to extend the "Unassigned" placeholder row.
$ws.Range("K3").Value = "Unassigned"
